$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/10/2025  Through  3/16/2025"

# --- Reference cells used to copy style+type for text<->number conversions ---
# I16 = stable numeric (style 14) reference
# H16 = stable numeric (style 15, percent) reference
# C22 = stable text "0" (style 13, shared string index 20) reference
# E14 = stable text "***.*" (style 13, shared string index 21) reference

# --- Cell value / type updates ---
$ws.Range("F16").Value2 = 2
$ws.Range("G16").Value2 = 2
$ws.Range("L16").Value2 = -58.823529411764
$ws.Range("C17").Value2 = 4
$ws.Range("D17").Value2 = 7
$ws.Range("E17").Value2 = -42.857142857142
$ws.Range("G17").Value2 = 20
$ws.Range("H17").Value2 = -30
$ws.Range("I17").Value2 = 48
$ws.Range("J17").Value2 = 48
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 2.127659574468
$ws.Range("I16").Copy($ws.Range("C18"))
$ws.Range("C18").Value2 = 1
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 2
$ws.Range("G18").Value2 = 4
$ws.Range("I18").Value2 = 10
$ws.Range("J18").Value2 = 12
$ws.Range("K18").Value2 = -16.666666666666
$ws.Range("L18").Value2 = -47.368421052631
$ws.Range("C19").Value2 = 1
$ws.Range("D19").Value2 = 7
$ws.Range("E19").Value2 = -85.714285714285
$ws.Range("G19").Value2 = 18
$ws.Range("H19").Value2 = 5.555555555555
$ws.Range("I19").Value2 = 41
$ws.Range("J19").Value2 = 65
$ws.Range("K19").Value2 = -36.923076923076
$ws.Range("L19").Value2 = -25.454545454545
$ws.Range("C20").Value2 = 4
$ws.Range("I16").Copy($ws.Range("D20"))
$ws.Range("D20").Value2 = 1
$ws.Range("H16").Copy($ws.Range("E20"))
$ws.Range("E20").Value2 = 300
$ws.Range("F20").Value2 = 12
$ws.Range("G20").Value2 = 7
$ws.Range("H20").Value2 = 71.428571428571
$ws.Range("I20").Value2 = 39
$ws.Range("J20").Value2 = 32
$ws.Range("K20").Value2 = 21.875
$ws.Range("L20").Value2 = 18.181818181818
$ws.Range("C21").Value2 = 10
$ws.Range("D21").Value2 = 16
$ws.Range("E21").Value2 = -37.5
$ws.Range("F21").Value2 = 50
$ws.Range("G21").Value2 = 52
$ws.Range("H21").Value2 = -3.846153846153
$ws.Range("I21").Value2 = 152
$ws.Range("J21").Value2 = 179
$ws.Range("K21").Value2 = -15.083798882681
$ws.Range("L21").Value2 = -11.627906976744
$ws.Range("C24").Value2 = 10
$ws.Range("D24").Value2 = 11
$ws.Range("E24").Value2 = -9.090909090909
$ws.Range("F24").Value2 = 33
$ws.Range("G24").Value2 = 49
$ws.Range("H24").Value2 = -32.653061224489
$ws.Range("I24").Value2 = 104
$ws.Range("J24").Value2 = 122
$ws.Range("K24").Value2 = -14.754098360655
$ws.Range("L24").Value2 = -17.460317460317
$ws.Range("C22").Copy($ws.Range("C25"))
$ws.Range("D25").Value2 = 4
$ws.Range("E25").Value2 = -100
$ws.Range("F25").Value2 = 5
$ws.Range("G25").Value2 = 9
$ws.Range("H25").Value2 = -44.444444444444
$ws.Range("J25").Value2 = 18
$ws.Range("K25").Value2 = 38.888888888888
$ws.Range("L25").Value2 = 8.695652173913
$ws.Range("C26").Value2 = 7
$ws.Range("D26").Value2 = 10
$ws.Range("E26").Value2 = -30
$ws.Range("F26").Value2 = 26
$ws.Range("G26").Value2 = 50
$ws.Range("H26").Value2 = -48
$ws.Range("I26").Value2 = 77
$ws.Range("J26").Value2 = 97
$ws.Range("K26").Value2 = -20.61855670103
$ws.Range("L26").Value2 = 40
$ws.Range("C22").Copy($ws.Range("C27"))
$ws.Range("F27").Value2 = 2
$ws.Range("G27").Value2 = 1
$ws.Range("H27").Value2 = 100
$ws.Range("L27").Value2 = 80
$ws.Range("C22").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("L28").Value2 = -60
$ws.Range("C22").Copy($ws.Range("F33"))
